$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price text stays text (matches original inlineStr cells)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '21.667.44'
$ws.Range("E2").Value = '  -1.55%  '
$ws.Range("D3").Value = '1.532.99'
$ws.Range("E3").Value = '  -1.40%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").Value = '288.31'
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("E7").Value = '  +3.26%  '
$ws.Range("D8").Value = '0.3158'
$ws.Range("E8").Value = '  -2.47%  '
$ws.Range("D9").Value = '42.26'
$ws.Range("E9").Value = '  +2.26%  '
$ws.Range("D10").Value = '0.07156'
$ws.Range("E10").Value = '  -2.19%  '
$ws.Range("D11").Value = '1.044'
$ws.Range("E11").Value = '  -6.84%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("D13").Value = '5.619'
$ws.Range("E13").Value = '  -1.69%  '
$ws.Range("D14").Value = '18.43'
$ws.Range("E14").Value = '  -4.83%  '
$ws.Range("D15").Value = '6.593'
$ws.Range("E15").Value = '  -3.18%  '
$ws.Range("D16").Value = '1.537.32'
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("E18").Value = '  -0.47%  '
$ws.Range("D19").Value = '82.91'
$ws.Range("E19").Value = '  -2.67%  '
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").Value = '6.085'
$ws.Range("D22").Value = '15.38'
$ws.Range("E22").Value = '  -3.43%  '
$ws.Range("D23").Value = '10.79'
$ws.Range("E23").Value = '  -5.71%  '
$ws.Range("D24").Value = '2.389'
$ws.Range("E24").Value = '  +4.21%  '
$ws.Range("D25").Value = '21.674.96'
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("D26").Value = '2.325'
$ws.Range("E26").Value = '  -7.89%  '
$ws.Range("D27").Value = '147.54'
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").Value = '18.27'
$ws.Range("E28").Value = '  -3.02%  '
$ws.Range("D29").Value = '4.844'
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").Value = '1.713.39'
$ws.Range("E30").Value = '  -1.07%  '
$ws.Range("D31").Value = '116.77'
$ws.Range("E31").Value = '  -3.02%  '
$ws.Range("D32").Value = '5.834'
$ws.Range("E32").Value = '  -0.71%  '
$ws.Range("D33").Value = '0.9514'
$ws.Range("E33").Value = '  -13.05%  '
$ws.Range("D34").Value = '0.08111'
$ws.Range("E34").Value = '  -0.32%  '
$ws.Range("D35").Value = '8.525'
$ws.Range("E35").Value = '  -8.13%  '
$ws.Range("D36").Value = '0.06050'
$ws.Range("D37").Value = '5.077'
$ws.Range("E37").Value = '  -3.18%  '
$ws.Range("D38").Value = '0.02201'
$ws.Range("E38").Value = '  -4.04%  '
$ws.Range("D39").Value = '1.446'
$ws.Range("E39").Value = '  -11.97%  '
$ws.Range("E40").Value = '  -4.13%  '
$ws.Range("D41").Value = '1.175'
$ws.Range("E41").Value = '  -3.52%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D43").Value = '10.80'
$ws.Range("E43").Value = '  -0.49%  '
$ws.Range("D44").Value = '0.5723'
$ws.Range("E44").Value = '  -3.49%  '
$ws.Range("D45").Value = '12.99'
$ws.Range("E45").Value = '  -3.79%  '
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Value = '0.5473'
$ws.Range("E47").Value = '  -4.62%  '
$ws.Range("D48").Value = '1.153'
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("D49").Value = '1.858'
$ws.Range("D50").Value = '115.17'
$ws.Range("E50").Value = '  -3.42%  '
$ws.Range("D51").Value = '0.06682'
$ws.Range("E51").Value = '  -2.68%  '
